# Applies the registration-sheet update described by the commit:
#   - A new registrant (AGNIVA BHATTACHARJEE) is appended as row 46.
#   - The stray formula on F45 (phone number wrongly entered as a
#     formula) is removed, leaving just the text value.
#
# Helper: writes $val into $addr as a genuine text cell (no leftover
# formula, no forced "stored as text" style) by routing the literal
# through a temporary ="..." formula and then collapsing it to a
# value via Copy/PasteSpecial(values). A straight `.Value = $val`
# assignment would let Excel "helpfully" reinterpret digit-strings
# like "2020" or "08420880979" as numbers (losing the leading zero)
# or dates like "2002-01-21" as date serials, which is not what the
# source workbook stores.
function Set-TextValue($addr, $val) {
    $escaped = $val -replace '"', '""'
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F45: drop the erroneous "=08420880979" formula, keep the text value ---
Set-TextValue "F45" "08420880979"

# --- Row 46: newly registered entry ---
Set-TextValue "A46" "BSS/9bf9433c0000"
Set-TextValue "B46" "AGNIVA"
Set-TextValue "C46" "BHATTACHARJEE"
Set-TextValue "D46" "something"
Set-TextValue "E46" "2020"
$ws.Range("F46").Formula = "=08420880979"
Set-TextValue "G46" "bhattacharjee.agniva.jobs@gmail.com"
Set-TextValue "H46" "2002-01-21"
Set-TextValue "I46" "IT"
Set-TextValue "J46" "Google"
Set-TextValue "N46" "aergty7u6543"
# K46, L46, M46 (Highest Qualification / Institution / Fondest Memory)
# are left blank, matching the empty cells in the source row.

$excel.CutCopyMode = $false

# Extend the "numbers stored as text" warning suppression to cover the
# new row, mirroring the widened dimension/ignoredErrors range.
$ws.Range("A1:N46").Errors.Item(-2146826246).Ignore = $true
